$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 4
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -2
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 3
